$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.052073
$ws.Range("H2").Value = 0.156219
$ws.Range("I2").Value = 0.03816600682360385
$ws.Range("J2").Value = 0.03816600682360385
$ws.Range("M2").Value = 28.31444233333334
$ws.Range("N2").Value = 84.94332700000001
$ws.Range("O2").Value = 0.2747173016130739
$ws.Range("P2").Value = 0.2747173016130739
$ws.Range("Q2").Value = 1.474417955623667
$ws.Range("R2").Value = 13.269761600613
$ws.Range("S2").Value = 0.01048486240792662
$ws.Range("T2").Value = 0.01048486240792662
$ws.Range("G3").Value = 0.052073
$ws.Range("H3").Value = 0.156219
$ws.Range("I3").Value = 0.03816600682360385
$ws.Range("J3").Value = 0.03816600682360385
$ws.Range("O3").Value = 0.2090339131726295
$ws.Range("P3").Value = 0.2090339131726295
$ws.Range("Q3").Value = 1.121892771610334
$ws.Range("R3").Value = 10.097034944493
$ws.Range("S3").Value = 0.007977989756511194
$ws.Range("T3").Value = 0.007977989756511194
$ws.Range("G4").Value = 0.052073
$ws.Range("H4").Value = 0.156219
$ws.Range("I4").Value = 0.03816600682360385
$ws.Range("J4").Value = 0.03816600682360385
$ws.Range("M4").Value = 5.413469333333334
$ws.Range("N4").Value = 16.240408
$ws.Range("O4").Value = 0.0525235026743817
$ws.Range("P4").Value = 0.0525235026743817
$ws.Range("Q4").Value = 0.2818955885946667
$ws.Range("R4").Value = 2.537060297352
$ws.Range("S4").Value = 0.002004612361470027
$ws.Range("T4").Value = 0.002004612361470027
$ws.Range("G5").Value = 0.052073
$ws.Range("H5").Value = 0.156219
$ws.Range("I5").Value = 0.03816600682360385
$ws.Range("J5").Value = 0.03816600682360385
$ws.Range("M5").Value = 47.79503400000001
$ws.Range("N5").Value = 143.385102
$ws.Range("O5").Value = 0.4637252825399149
$ws.Range("P5").Value = 0.4637252825399149
$ws.Range("Q5").Value = 2.488830805482
$ws.Range("R5").Value = 22.399477249338
$ws.Range("S5").Value = 0.01769854229769601
$ws.Range("T5").Value = 0.01769854229769601
$ws.Range("I6").Value = 0.5197685398391702
$ws.Range("J6").Value = 0.5197685398391702
$ws.Range("M6").Value = 28.31444233333334
$ws.Range("N6").Value = 84.94332700000001
$ws.Range("O6").Value = 0.2747173016130739
$ws.Range("P6").Value = 0.2747173016130739
$ws.Range("Q6").Value = 20.07954543028622
$ws.Range("R6").Value = 180.715908872576
$ws.Range("S6").Value = 0.1427894107279843
$ws.Range("T6").Value = 0.1427894107279843
$ws.Range("I7").Value = 0.5197685398391702
$ws.Range("J7").Value = 0.5197685398391702
$ws.Range("O7").Value = 0.2090339131726295
$ws.Range("P7").Value = 0.2090339131726295
$ws.Range("S7").Value = 0.1086492518266055
$ws.Range("T7").Value = 0.1086492518266055
$ws.Range("I8").Value = 0.5197685398391702
$ws.Range("J8").Value = 0.5197685398391702
$ws.Range("M8").Value = 5.413469333333334
$ws.Range("N8").Value = 16.240408
$ws.Range("O8").Value = 0.0525235026743817
$ws.Range("P8").Value = 0.0525235026743817
$ws.Range("Q8").Value = 3.83903034834489
$ws.Range("R8").Value = 34.551273135104
$ws.Range("S8").Value = 0.02730006429230213
$ws.Range("T8").Value = 0.02730006429230213
$ws.Range("I9").Value = 0.5197685398391702
$ws.Range("J9").Value = 0.5197685398391702
$ws.Range("M9").Value = 47.79503400000001
$ws.Range("N9").Value = 143.385102
$ws.Range("O9").Value = 0.4637252825399149
$ws.Range("P9").Value = 0.4637252825399149
$ws.Range("Q9").Value = 33.89445376486401
$ws.Range("R9").Value = 305.050083883776
$ws.Range("S9").Value = 0.2410298129922782
$ws.Range("T9").Value = 0.2410298129922782
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.299804
$ws.Range("H10").Value = 0.8994119999999999
$ws.Range("I10").Value = 0.2197361686429383
$ws.Range("J10").Value = 0.2197361686429384
$ws.Range("M10").Value = 28.31444233333334
$ws.Range("N10").Value = 84.94332700000001
$ws.Range("O10").Value = 0.2747173016130739
$ws.Range("P10").Value = 0.2747173016130739
$ws.Range("Q10").Value = 8.488783069302666
$ws.Range("R10").Value = 76.399047623724
$ws.Range("S10").Value = 0.06036532731638337
$ws.Range("T10").Value = 0.06036532731638337
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.299804
$ws.Range("H11").Value = 0.8994119999999999
$ws.Range("I11").Value = 0.2197361686429383
$ws.Range("J11").Value = 0.2197361686429384
$ws.Range("O11").Value = 0.2090339131726295
$ws.Range("P11").Value = 0.2090339131726295
$ws.Range("Q11").Value = 6.459161955329333
$ws.Range("R11").Value = 58.13245759796399
$ws.Range("S11").Value = 0.04593231119699424
$ws.Range("T11").Value = 0.04593231119699425
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.299804
$ws.Range("H12").Value = 0.8994119999999999
$ws.Range("I12").Value = 0.2197361686429383
$ws.Range("J12").Value = 0.2197361686429384
$ws.Range("M12").Value = 5.413469333333334
$ws.Range("N12").Value = 16.240408
$ws.Range("O12").Value = 0.0525235026743817
$ws.Range("P12").Value = 0.0525235026743817
$ws.Range("Q12").Value = 1.622979760010667
$ws.Range("R12").Value = 14.606817840096
$ws.Range("S12").Value = 0.01154131324137576
$ws.Range("T12").Value = 0.01154131324137576
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.299804
$ws.Range("H13").Value = 0.8994119999999999
$ws.Range("I13").Value = 0.2197361686429383
$ws.Range("J13").Value = 0.2197361686429384
$ws.Range("M13").Value = 47.79503400000001
$ws.Range("N13").Value = 143.385102
$ws.Range("O13").Value = 0.4637252825399149
$ws.Range("P13").Value = 0.4637252825399149
$ws.Range("Q13").Value = 14.329142373336
$ws.Range("R13").Value = 128.962281360024
$ws.Range("S13").Value = 0.101897216888185
$ws.Range("T13").Value = 0.101897216888185
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.303342
$ws.Range("H14").Value = 0.910026
$ws.Range("I14").Value = 0.2223292846942876
$ws.Range("J14").Value = 0.2223292846942876
$ws.Range("M14").Value = 28.31444233333334
$ws.Range("N14").Value = 84.94332700000001
$ws.Range("O14").Value = 0.2747173016130739
$ws.Range("P14").Value = 0.2747173016130739
$ws.Range("Q14").Value = 8.588959566278001
$ws.Range("R14").Value = 77.300636096502
$ws.Range("S14").Value = 0.0610777011607796
$ws.Range("T14").Value = 0.0610777011607796
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.303342
$ws.Range("H15").Value = 0.910026
$ws.Range("I15").Value = 0.2223292846942876
$ws.Range("J15").Value = 0.2223292846942876
$ws.Range("O15").Value = 0.2090339131726295
$ws.Range("P15").Value = 0.2090339131726295
$ws.Range("Q15").Value = 6.535386805558001
$ws.Range("R15").Value = 58.818481250022
$ws.Range("S15").Value = 0.04647436039251854
$ws.Range("T15").Value = 0.04647436039251855
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.303342
$ws.Range("H16").Value = 0.910026
$ws.Range("I16").Value = 0.2223292846942876
$ws.Range("J16").Value = 0.2223292846942876
$ws.Range("M16").Value = 5.413469333333334
$ws.Range("N16").Value = 16.240408
$ws.Range("O16").Value = 0.0525235026743817
$ws.Range("P16").Value = 0.0525235026743817
$ws.Range("Q16").Value = 1.642132614512
$ws.Range("R16").Value = 14.779193530608
$ws.Range("S16").Value = 0.01167751277923379
$ws.Range("T16").Value = 0.01167751277923379
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.303342
$ws.Range("H17").Value = 0.910026
$ws.Range("I17").Value = 0.2223292846942876
$ws.Range("J17").Value = 0.2223292846942876
$ws.Range("M17").Value = 47.79503400000001
$ws.Range("N17").Value = 143.385102
$ws.Range("O17").Value = 0.4637252825399149
$ws.Range("P17").Value = 0.4637252825399149
$ws.Range("Q17").Value = 14.498241203628
$ws.Range("R17").Value = 130.484170832652
$ws.Range("S17").Value = 0.1030997103617557
$ws.Range("T17").Value = 0.1030997103617557

Write-Output "Applied 190 cell updates"
